$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics1")

$ws.Range("H1").Value = "negative effectors"
$ws.Range("I1").Value = "positive effectors"

$ws.Range("L1").Value = "mechanism_refs_type"
$ws.Range("M1").Value = "mechanism_refs"
$ws.Range("N1").Value = "inhibitors_refs_type"
$ws.Range("O1").Value = "inhibitors_refs"
$ws.Range("P1").Value = "activators_refs_type"
$ws.Range("Q1").Value = "activators_refs"
$ws.Range("R1").Value = "negative_effectors_refs_type"
$ws.Range("S1").Value = "negative_effectors_refs"
$ws.Range("T1").Value = "positive_effectors_refs_type"
$ws.Range("U1").Value = "positive_effectors_refs"
$ws.Range("V1").Value = "subunits_refs_type"
$ws.Range("W1").Value = "subunits_refs"
